# "Generate Report for Archive"
#
# 1. Status text "Ready for handoff" -> "In Translation" everywhere it
#    appears: Overview!E2:F3 (the zh-cn / de-de status columns) and the
#    Status column (C) on both the zh-cn and de-de detail sheets.
# 2. Narrow the (now shorter) status columns: Overview columns E & F, and
#    column C on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# New narrower width (characters) for the status columns - matches the
# ~13.41-character column the authoring app wrote for the shorter status
# text (this COM layer quantizes ColumnWidth to the nearest 1/6 character,
# so 12.5 is the input that lands closest to that stored width).
$statusColWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
